$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 3 new rows before row 97, copying formatting from row 96 (the row above)
$ws.Rows.Item(97).Insert()
$ws.Rows.Item(97).Insert()
$ws.Rows.Item(97).Insert()
$ws.Range("A96:K96").Copy()
$ws.Range("A97:K99").PasteSpecial(-4122)

# Resize Table1 to include the new rows
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K135"))

# Restore calculated-column formula on the new rows' EARNED' (G) column
$ws.Range("G97").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G98").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G99").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G133").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G134").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("G135").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Row 96: new SP(1-0-0) entry earning 1.25, remark BDAY 4/3/23
$ws.Range("B96").Value = "SP(1-0-0)"
$ws.Range("C96").Value = 1.25
$ws.Range("K96").Value = "BDAY 4/3/23"

# Row 97: SP(1-0-0) entry, remark ANNIV 4/4/23
$ws.Range("B97").Value = "SP(1-0-0)"
$ws.Range("K97").Value = "ANNIV 4/4/23"

# Row 98: SP(1-0-0) entry, remark FILIAL 4/5/23
$ws.Range("B98").Value = "SP(1-0-0)"
$ws.Range("K98").Value = "FILIAL 4/5/23"

# Row 99: VL(2-0-0) entry, 2 days absence w/o pay, remark 4/11,12/2023
$ws.Range("B99").Value = "VL(2-0-0)"
$ws.Range("D99").Value = 2
$ws.Range("K99").Value = "4/11,12/2023"
